$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "Packages"
$ws.Range("C3").Value = "Packages"

$ws.Rows.Item(2).RowHeight = 42.75
$ws.Rows.Item(3).RowHeight = 42.75

$ws.Range("C3").Select()
